# detail Info almost finish
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Select 부분" row: mark as done (X -> O)
$ws.Range("B9").Value = "O"

# Move the active selection to H8, as left by the author after the edit
$ws.Range("H8").Select()
